$wb = $excel.ActiveWorkbook

# Update "Forecast Comparison" sheet: MyForecast value for W2 (D3) 13 -> 12
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsForecast.Range("D3").Value = 12

# Update "Summary" sheet: Total Forecast figures (stored as text strings)
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "306"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "140"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "57"
